$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header row: add x, y columns
$ws.Range("C1").Value = "x"
$ws.Range("D1").Value = "y"

# Set column widths for the new columns (C and D) to match A and B (raw width 15)
$ws.Columns("C").ColumnWidth = 14.1666666667
$ws.Columns("D").ColumnWidth = 14.1666666667

# Data values: fid, Type, x, y
$data = @(
    @(1, "Koud", 96602, 437272),
    @(2, "Koud", 96588, 437350),
    @(3, "Koud", 96668, 437322),
    @(4, "Koud", 96527, 437289),
    @(6, "Warm", 96776, 437219),
    @(7, "Warm", 96840, 437279),
    @(8, "Warm", 96870, 437196),
    @(9, "Warm", 96927, 437260)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $row++
}

# Remove the extra 10th row (old data had 9 rows of data + header = 10 rows total)
$ws.Rows(10).Delete()
